# TimeLog_ConnorPeper.xlsx - "Add files via upload" edit
#
# Semantic changes applied (derived from the OOXML diff):
#   1. Row 15 (week 9) "Hours" value changes from 2 to 4.
#   2. Row 15 (week 9) "Activities" note text is extended with additional
#      status updates (Task 16 / Task 26 complete, server setup notes, etc.).
#   3. The "Hours" column (E) is reformatted to show two decimal places
#      (number format 0.00) for every row that already holds data.
#
# (Cosmetic, machine/session specific artifacts such as the scroll
# position, Excel build metadata, and the localized "Normal"/"Standard"
# cell-style name are not part of the authored content and are left to
# the runtime's own save behaviour.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1 & 2: update the Week 9 row (hours worked + activity log text)
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = "Task 4 & 6: Completed backend and database implementation of message updating. Task 16: Laid out the groundwork for updating. Downloaded all necessary software to develop stuff without my own computer. All of this is on the Server. Reformatted some code to a different response body we might want to use. Not being used yet. Task 16: Complete. Task 26: Complete"

# 3: apply a 2 decimal place number format to the populated Hours cells
$ws.Range("E3").NumberFormat = "0.00"
$ws.Range("E4").NumberFormat = "0.00"
$ws.Range("E6:E22").NumberFormat = "0.00"

Write-Output "edit complete"
